$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Pearson logo (footers): image2.png -> image1.png
for ($i = 1; $i -le 2; $i++) {
    $f = $sec.Footers.Item($i)
    if ($f.Exists -and $f.Range.InlineShapes.Count -gt 0) {
        for ($j = 1; $j -le $f.Range.InlineShapes.Count; $j++) {
            $shp = $f.Range.InlineShapes.Item($j)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}

# BTec logo (headers): image1.jpg -> image2.jpg
for ($i = 1; $i -le 2; $i++) {
    $h = $sec.Headers.Item($i)
    if ($h.Exists -and $h.Range.InlineShapes.Count -gt 0) {
        for ($j = 1; $j -le $h.Range.InlineShapes.Count; $j++) {
            $shp = $h.Range.InlineShapes.Item($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }
}
